# Update to latest business count data (2024)
# Source: ONS UK Business Counts (rows 8 & 9 - "Enterprises by employment size
# band" / "Enterprises by employment industry"). The "Latest period" (col C)
# and "Next period" (col D) release dates move on one year.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: Enterprises by employment size band
$ws.Cells.Item(8, 3).Value = "Mar 2024 (17/10/24)"
$ws.Cells.Item(8, 4).Value = "Mar 2025 (30/12/25)"

# Row 9: Enterprises by employment industry
$ws.Cells.Item(9, 3).Value = "Mar 2024 (17/10/24)"
$ws.Cells.Item(9, 4).Value = "Mar 2025 (30/12/25)"

# Reflect the author's saved view state - selection moved to the updated
# "Latest period" cell for row 9 (the workbook was left scrolled so row 3
# is at the top; this COM layer does not expose a settable scroll-position
# API, so only the selection itself is reproduced here).
$ws.Range("C9").Select()
